{"js": "const replacements = [\n  [\"452\u00f79=50, 2\", \"129\u00f72=64, 1\"],\n  [\"447\u00f76=74, 3\", \"977\u00f76=162, 5\"],\n  [\"963\u00f75=192, 3\", \"281\u00f78=35, 1\"],\n  [\"578\u00f76=96, 2\", \"500\u00f75=100, 0\"],\n  [\"745\u00f72=372, 1\", \"339\u00f73=113, 0\"],\n  [\"639\u00f77=91, 2\", \"428\u00f75=85, 3\"],\n  [\"728\u00f74=182, 0\", \"853\u00f77=121, 6\"],\n  [\"806\u00f73=268, 2\", \"987\u00f72=493, 1\"],\n  [\"957\u00f78=119, 5\", \"883\u00f76=147, 1\"],\n  [\"747\u00f79=83, 0\", \"675\u00f75=135, 0\"],\n  [\"231\u00f74=57, 3\", \"828\u00f77=118, 2\"],\n  [\"353\u00f72=176, 1\", \"816\u00f77=116, 4\"],\n  [\"998\u00f78=124, 6\", \"728\u00f77=104, 0\"],\n  [\"995\u00f76=165, 5\", \"662\u00f76=110, 2\"],\n  [\"583\u00f76=97, 1\", \"358\u00f79=39, 7\"],\n  [\"455\u00f74=113, 3\", \"619\u00f73=206, 1\"],\n  [\"950\u00f73=316, 2\", \"267\u00f72=133, 1\"],\n  [\"745\u00f76=124, 1\", \"929\u00f78=116, 1\"],\n  [\"253\u00f79=28, 1\", \"668\u00f72=334, 0\"],\n  [\"755\u00f77=107, 6\", \"250\u00f76=41, 4\"],\n  [\"908\u00f72=454, 0\", \"172\u00f76=28, 4\"],\n  [\"690\u00f74=172, 2\", \"857\u00f74=214, 1\"],\n  [\"103\u00f79=11, 4\", \"420\u00f79=46, 6\"],\n  [\"918\u00f74=229, 2\", \"785\u00f75=157, 0\"],\n  [\"860\u00f75=172, 0\", \"593\u00f73=197, 2\"],\n];\n\nconst body = context.document.body;\nfor (const [oldText, newText] of replacements) {\n  const results = body.search(oldText, { matchCase: true, matchWholeWord: false });\n  results.load(\"items\");\n  await context.sync();\n  for (const item of results.items) {\n    item.insertText(newText, \"Replace\");\n  }\n  await context.sync();\n}\n", "ps1": "$d = $word.ActiveDocument\n\n$pairs = @(\n  @(\"452\u00f79=50, 2\", \"129\u00f72=64, 1\"),\n  @(\"447\u00f76=74, 3\", \"977\u00f76=162, 5\"),\n  @(\"963\u00f75=192, 3\", \"281\u00f78=35, 1\"),\n  @(\"578\u00f76=96, 2\", \"500\u00f75=100, 0\"),\n  @(\"745\u00f72=372, 1\", \"339\u00f73=113, 0\"),\n  @(\"639\u00f77=91, 2\", \"428\u00f75=85, 3\"),\n  @(\"728\u00f74=182, 0\", \"853\u00f77=121, 6\"),\n  @(\"806\u00f73=268, 2\", \"987\u00f72=493, 1\"),\n  @(\"957\u00f78=119, 5\", \"883\u00f76=147, 1\"),\n  @(\"747\u00f79=83, 0\", \"675\u00f75=135, 0\"),\n  @(\"231\u00f74=57, 3\", \"828\u00f77=118, 2\"),\n  @(\"353\u00f72=176, 1\", \"816\u00f77=116, 4\"),\n  @(\"998\u00f78=124, 6\", \"728\u00f77=104, 0\"),\n  @(\"995\u00f76=165, 5\", \"662\u00f76=110, 2\"),\n  @(\"583\u00f76=97, 1\", \"358\u00f79=39, 7\"),\n  @(\"455\u00f74=113, 3\", \"619\u00f73=206, 1\"),\n  @(\"950\u00f73=316, 2\", \"267\u00f72=133, 1\"),\n  @(\"745\u00f76=124, 1\", \"929\u00f78=116, 1\"),\n  @(\"253\u00f79=28, 1\", \"668\u00f72=334, 0\"),\n  @(\"755\u00f77=107, 6\", \"250\u00f76=41, 4\"),\n  @(\"908\u00f72=454, 0\", \"172\u00f76=28, 4\"),\n  @(\"690\u00f74=172, 2\", \"857\u00f74=214, 1\"),\n  @(\"103\u00f79=11, 4\", \"420\u00f79=46, 6\"),\n  @(\"918\u00f74=229, 2\", \"785\u00f75=157, 0\"),\n  @(\"860\u00f75=172, 0\", \"593\u00f73=197, 2\"),\n)\n\nforeach ($pair in $pairs) {\n  $find = $d.Content.Find\n  $find.Execute($pair[0], $false, $false, $false, $false, $false, $true, 1, $false, $pair[1], 2)\n}\n"}
